$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A186").Value = "IMX-USD"
$ws.Range("A187").Value = "TAO-USD"
$ws.Range("A188").Value = "GRT-USD"
